# Regenerate the handback-status report: the handoff/handback run that
# produced this workbook picked up new source GUID-based file names and new
# timestamps for both the "Correspond Handoff" (xliff generation) and the
# "Correspond Handback" (translation ingestion) events. Update every sheet
# (Overview, zh-cn, de-de) with the refreshed file names / dates, and keep
# each hyperlink's display text in sync with its cell text.

$wb = $excel.ActiveWorkbook

$newName1 = "ef46a31a-eafe-475f-b8e9-a54f33cbe7d2"
$newName2 = "ffffdd49c18f-f6e3-4132-8fdd-56e83a1fbfd8"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = ($newName1 + ".md")
$ws.Range("B2").Value = ("e2e\" + $newName1 + ".md")
$ws.Range("G2").Value = "2016-09-05 09:26:30"

$ws.Range("A3").Value = ($newName2 + ".md")
$ws.Range("B3").Value = ("e2e\" + $newName2 + ".md")
$ws.Range("G3").Value = "2016-09-05 09:26:30"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$2') {
        $h.TextToDisplay = ("e2e\" + $newName1 + ".md")
    } elseif ($addr -eq '$B$3') {
        $h.TextToDisplay = ("e2e\" + $newName2 + ".md")
    }
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = ($newName1 + ".md")
$ws.Range("I2").Value = ($newName1 + ".md")
$ws.Range("G2").Value = ($newName1 + ".f40c5855785e29b39e864cd0988ce70e08c92d7f.zh-cn.xlf")
$ws.Range("H2").Value = "2016-09-05 09:26:25"
$ws.Range("J2").Value = ($newName1 + ".f40c5855785e29b39e864cd0988ce70e08c92d7f.zh-cn.xlf")
$ws.Range("K2").Value = "2016-09-05 09:26:55"

$ws.Range("A3").Value = ($newName2 + ".md")
$ws.Range("I3").Value = ($newName2 + ".md")
$ws.Range("G3").Value = ($newName1 + ".f40c5855785e29b39e864cd0988ce70e08c92d7f.zh-cn.xlf")
$ws.Range("H3").Value = "2016-09-05 09:26:25"
$ws.Range("J3").Value = ($newName1 + ".f40c5855785e29b39e864cd0988ce70e08c92d7f.zh-cn.xlf")
$ws.Range("K3").Value = "2016-09-05 09:26:55"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = ($newName1 + ".md")
    } elseif ($addr -eq '$I$2') {
        $h.TextToDisplay = ($newName1 + ".md")
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = ($newName2 + ".md")
    } elseif ($addr -eq '$I$3') {
        $h.TextToDisplay = ($newName2 + ".md")
    }
}

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = ($newName1 + ".md")
$ws.Range("I2").Value = ($newName1 + ".md")
$ws.Range("G2").Value = ($newName1 + ".f40c5855785e29b39e864cd0988ce70e08c92d7f.de-de.xlf")
$ws.Range("H2").Value = "2016-09-05 09:26:30"
$ws.Range("J2").Value = ($newName1 + ".f40c5855785e29b39e864cd0988ce70e08c92d7f.de-de.xlf")
$ws.Range("K2").Value = "2016-09-05 09:27:09"

$ws.Range("A3").Value = ($newName2 + ".md")
$ws.Range("I3").Value = ($newName2 + ".md")
$ws.Range("G3").Value = ($newName1 + ".f40c5855785e29b39e864cd0988ce70e08c92d7f.de-de.xlf")
$ws.Range("H3").Value = "2016-09-05 09:26:30"
$ws.Range("J3").Value = ($newName1 + ".f40c5855785e29b39e864cd0988ce70e08c92d7f.de-de.xlf")
$ws.Range("K3").Value = "2016-09-05 09:27:09"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = ($newName1 + ".md")
    } elseif ($addr -eq '$I$2') {
        $h.TextToDisplay = ($newName1 + ".md")
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = ($newName2 + ".md")
    } elseif ($addr -eq '$I$3') {
        $h.TextToDisplay = ($newName2 + ".md")
    }
}
